# "improve codes and add more testing for file api"
#
# - Strips the bold/bordered/centered header style from row 1 (back to the
#   default "Normal" style).
# - Row 2's numeric-looking cells (id, purchase_amount, shipping_zip,
#   payment_amount, payment_final_capture, paypal_order_id) are rewritten
#   as text while keeping the same displayed value.
# - Four new "testing" columns (U:X) are appended to row 2 with the text
#   value "s".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:T1): remove the custom font/border/alignment style.
$ws.Range("A1:T1").Style = "Normal"

# Row 2 cells that need to flip from numeric storage to text storage,
# without changing what they display. A leading apostrophe forces Excel to
# store the cell as text; re-applying the "Normal" style afterwards clears
# the transient quote-prefix formatting it introduces, leaving a plain
# text cell that uses the default style (matches the rest of the row).
$textCells = @(
    @{ Addr = "A2"; Value = "1" },
    @{ Addr = "D2"; Value = "1.77" },
    @{ Addr = "L2"; Value = "95131" },
    @{ Addr = "Q2"; Value = "1.77" },
    @{ Addr = "R2"; Value = "1" },
    @{ Addr = "T2"; Value = "1" }
)
foreach ($cell in $textCells) {
    $rng = $ws.Range($cell.Addr)
    $rng.Value = "'" + $cell.Value
    $rng.Style = "Normal"
}

# New testing columns appended after the existing data (U2:X2), all "s".
foreach ($addr in @("U2", "V2", "W2", "X2")) {
    $ws.Range($addr).Value = "s"
}
